# ---------------------------------------------------------------------------
# B6-PowerPoint.pptx edit replay
#
# 1) Three tables (on the slides that used to hold the "Table_0" custom
#    table style) get re-styled to the built-in table style
#    {5777DE4F-A758-4B9E-8F0F-BA99AF33C565}.
# 2) The deck's theme palette is reset back to the stock "Office" colors
#    (it had been carrying the "Integral" / Red Violet palette).
# ---------------------------------------------------------------------------

$p = $ppt.ActivePresentation

# --- 1) Re-style the three tables -------------------------------------------------
$newStyleId = "{5777DE4F-A758-4B9E-8F0F-BA99AF33C565}"

for ($slideIdx = 14; $slideIdx -le 16; $slideIdx++) {
    $slide = $p.Slides.Item($slideIdx)
    for ($shapeIdx = 1; $shapeIdx -le $slide.Shapes.Count; $shapeIdx++) {
        $shape = $slide.Shapes.Item($shapeIdx)
        if ($shape.HasTable) {
            $shape.Table.ApplyStyle($newStyleId)
        }
    }
}

# --- 2) Restore the stock "Office" theme colors -----------------------------------
$cs = $p.SlideMaster.ColorScheme
$cs.Item(3).RGB  = 6968388    # dk2       44546A
$cs.Item(4).RGB  = 15132391   # lt2       E7E6E6
$cs.Item(5).RGB  = 13998939   # accent1   5B9BD5
$cs.Item(6).RGB  = 3243501    # accent2   ED7D31
$cs.Item(7).RGB  = 10855845   # accent3   A5A5A5
$cs.Item(8).RGB  = 49407      # accent4   FFC000
$cs.Item(9).RGB  = 12874308   # accent5   4472C4
$cs.Item(10).RGB = 4697456    # accent6   70AD47
$cs.Item(11).RGB = 12673797   # hlink     0563C1
$cs.Item(12).RGB = 7491477    # folHlink  954F72
